$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the diff.
# Values are set in (row, new-value) pairs; format as plain numeric literals
# so PowerShell/Excel COM stores them as numbers, not strings.

# Row 97
$ws.Range("F97").Value = 52
$ws.Range("G97").Value = 4941.04
# Row 107
$ws.Range("F107").Value = 40
$ws.Range("G107").Value = 2812.8
# Row 111
$ws.Range("F111").Value = 237
$ws.Range("G111").Value = 15096.9
# Row 123
$ws.Range("F123").Value = 38
$ws.Range("G123").Value = 7664.98
# Row 133
$ws.Range("B133").Value = 202940.34
# Row 164
$ws.Range("F164").Value = 21
$ws.Range("G164").Value = 1039.08
# Row 165
$ws.Range("F165").Value = 86
$ws.Range("G165").Value = 4255.28
# Row 173
$ws.Range("F173").Value = 95
$ws.Range("G173").Value = 4229.4
# Row 176
$ws.Range("B176").Value = 13683.85
# Row 207
$ws.Range("F207").Value = 30
$ws.Range("G207").Value = 1184.7
# Row 214
$ws.Range("F214").Value = 7
$ws.Range("G214").Value = 897.96
# Row 216
$ws.Range("B216").Value = 9360.860000000001
# Row 314
$ws.Range("B314").Value = 57077
$ws.Range("D314").Value = 93.08
$ws.Range("E314").Value = 111.2
$ws.Range("F314").Value = 1
$ws.Range("G314").Value = 93.08
# Row 315
$ws.Range("B315").Value = 61610
$ws.Range("D315").Value = 102.71
$ws.Range("E315").Value = 122.71
$ws.Range("F315").Value = 83
$ws.Range("G315").Value = 8524.93
# Row 316
$ws.Range("F316").Value = 13
$ws.Range("G316").Value = 3467.49
# Row 324
$ws.Range("F324").Value = 164
$ws.Range("G324").Value = 18733.72
# Row 327
$ws.Range("F327").Value = 55
$ws.Range("G327").Value = 7538.85
# Row 333
$ws.Range("F333").Value = 127
$ws.Range("G333").Value = 6517.64
# Row 342
$ws.Range("F342").Value = 121
$ws.Range("G342").Value = 6307.73
# Row 344
$ws.Range("F344").Value = 124
$ws.Range("G344").Value = 15997.24
# Row 354
$ws.Range("F354").Value = 28
$ws.Range("G354").Value = 2830.8
# Row 356
$ws.Range("F356").Value = 6
$ws.Range("G356").Value = 1191.54
# Row 369
$ws.Range("F369").Value = 2
$ws.Range("G369").Value = 1051.9
# Row 370
$ws.Range("F370").Value = 13
$ws.Range("G370").Value = 2597.01
# Row 380
$ws.Range("B380").Value = 243561.39
# Row 443
$ws.Range("F443").Value = 30
$ws.Range("G443").Value = 1116.3
# Row 447
$ws.Range("B447").Value = 38032.49
# Row 453
$ws.Range("F453").Value = 117
$ws.Range("G453").Value = 11495.25
# Row 459
$ws.Range("F459").Value = 69
$ws.Range("G459").Value = 3270.6
# Row 462
$ws.Range("F462").Value = 37
$ws.Range("G462").Value = 6898.28
# Row 473
$ws.Range("B473").Value = 134668.48
# Row 493
$ws.Range("F493").Value = 542
$ws.Range("G493").Value = 6943.02
# Row 505
$ws.Range("F505").Value = 403
$ws.Range("G505").Value = 5299.45
# Row 506
$ws.Range("F506").Value = 319
$ws.Range("G506").Value = 8389.700000000001
# Row 507
$ws.Range("F507").Value = 248
$ws.Range("G507").Value = 4074.64
# Row 509
$ws.Range("B509").Value = 93022.07000000001
# Row 511
$ws.Range("F511").Value = 10
$ws.Range("G511").Value = 369.7
# Row 516
$ws.Range("B516").Value = 6249.51
# Row 558
$ws.Range("F558").Value = 589
$ws.Range("G558").Value = 11691.65
# Row 560
$ws.Range("F560").Value = 291
$ws.Range("G560").Value = 4810.23
# Row 563
$ws.Range("B563").Value = 35977.59
# Row 619
$ws.Range("F619").Value = 20
$ws.Range("G619").Value = 2077.4
# Row 622
$ws.Range("F622").Value = 21
$ws.Range("G622").Value = 1054.83
# Row 640
$ws.Range("B640").Value = 205933.12
# Row 642
$ws.Range("F642").Value = 104
$ws.Range("G642").Value = 13577.2
# Row 646
$ws.Range("F646").Value = 2
$ws.Range("G646").Value = 54.4
# Row 649
$ws.Range("B649").Value = 52779.98
# Row 668
$ws.Range("F668").Value = 3
$ws.Range("G668").Value = 99.33
# Row 671
$ws.Range("F671").Value = 144
$ws.Range("G671").Value = 6217.92
# Row 673
$ws.Range("F673").Value = 57
$ws.Range("G673").Value = 2461.26
# Row 675
$ws.Range("F675").Value = 152
$ws.Range("G675").Value = 6563.36
# Row 677
$ws.Range("B677").Value = 19720.8
# Row 680
$ws.Range("F680").Value = 20
$ws.Range("G680").Value = 1597.6
# Row 681
$ws.Range("F681").Value = 9
$ws.Range("G681").Value = 738.36
# Row 682
$ws.Range("F682").Value = 21
$ws.Range("G682").Value = 1904.28
# Row 683
$ws.Range("F683").Value = 16
$ws.Range("G683").Value = 4974.08
# Row 685
$ws.Range("F685").Value = 19
$ws.Range("G685").Value = 1271.67
# Row 688
$ws.Range("F688").Value = 16
$ws.Range("G688").Value = 1520
# Row 693
$ws.Range("B693").Value = 29271.17
# Row 706
$ws.Range("F706").Value = 4
$ws.Range("G706").Value = 3786.84
# Row 714
$ws.Range("B714").Value = 84267.08
# Row 720
$ws.Range("F720").Value = 4
$ws.Range("G720").Value = 342
# Row 722
$ws.Range("B722").Value = 342
# Row 753
$ws.Range("F753").Value = 207
$ws.Range("G753").Value = 16882.92
# Row 756
$ws.Range("F756").Value = 235
$ws.Range("G756").Value = 30667.5
# Row 758
$ws.Range("F758").Value = 0
$ws.Range("G758").Value = 0
# Row 759
$ws.Range("F759").Value = 25
$ws.Range("G759").Value = 2788.5
# Row 761
$ws.Range("F761").Value = 80
$ws.Range("G761").Value = 1737.6
# Row 762
$ws.Range("F762").Value = 362
$ws.Range("G762").Value = 13495.36
# Row 769
$ws.Range("F769").Value = 446
$ws.Range("G769").Value = 60214.46
# Row 770
$ws.Range("F770").Value = 14
$ws.Range("G770").Value = 523.88
# Row 771
$ws.Range("F771").Value = 520
$ws.Range("G771").Value = 62769.2
# Row 772
$ws.Range("F772").Value = 41
$ws.Range("G772").Value = 4949.11
# Row 773
$ws.Range("B773").Value = 229797.45
# Row 850
$ws.Range("F850").Value = 520
$ws.Range("G850").Value = 15719.6
# Row 851
$ws.Range("F851").Value = 3047
$ws.Range("G851").Value = 496996.17
# Row 853
$ws.Range("F853").Value = 211
$ws.Range("G853").Value = 30521.15
# Row 854
$ws.Range("F854").Value = 111
$ws.Range("G854").Value = 4233.54
# Row 858
$ws.Range("F858").Value = 104
$ws.Range("G858").Value = 13372.32
# Row 859
$ws.Range("B859").Value = 603010.86
# Row 865
$ws.Range("B865").Value = 3361735.85
# Row 866
$ws.Range("B866").Value = 3361735.85
